# Insert one new data row for "Berenjena" at Vega Modelo de Temuco.
# The new row lands at sheet row 281 (pushing the old rows 281-318 down
# to 282-319), matching a weekly data refresh that prepends the newest
# observation to this variety's block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 281:318 down to 282:319, leaving a blank row 281.
$ws.Rows.Item(281).Insert()

# Populate the new row. Columns that are constant across this whole
# "Berenjena" block (A, B, C, E, F, G, H, I, N, Q, R) are copied from the
# neighboring row; the columns that actually vary (D, J, K, L, M, O, P)
# get the new observation's values.
$ws.Cells.Item(281, 1).Value = 10
$ws.Cells.Item(281, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(281, 3).Value = "La Araucanía"
$ws.Cells.Item(281, 4).Value = 44776
$ws.Cells.Item(281, 5).Value = 9
$ws.Cells.Item(281, 6).Value = 100112001
$ws.Cells.Item(281, 7).Value = "Berenjena"
$ws.Cells.Item(281, 8).Value = "Sin especificar"
$ws.Cells.Item(281, 9).Value = "Primera"
$ws.Cells.Item(281, 10).Value = 50
$ws.Cells.Item(281, 11).Value = 15000
$ws.Cells.Item(281, 12).Value = 15000
$ws.Cells.Item(281, 13).Value = 15000
$ws.Cells.Item(281, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(281, 15).Value = "Región del Maule"
$ws.Cells.Item(281, 16).Value = 250
$ws.Cells.Item(281, 17).Value = 60
$ws.Cells.Item(281, 18).Value = "Hortaliza"
